$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.4
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 9
$ws.Range("K2").Value = 2.2
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 1.33
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.36
$ws.Range("S2").Value = 2.1
$ws.Range("T2").Value = 1.7
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 1.38
$ws.Range("X2").Value = 1.25
$ws.Range("AA2").Value = 2.5
$ws.Range("AB2").Value = 1.5
$ws.Range("AE2").Value = 9.5
$ws.Range("AG2").Value = 15
$ws.Range("AM2").Value = 17
$ws.Range("AR2").Value = 81

# Row 3 updates
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2.4
$ws.Range("T3").Value = 1.53
$ws.Range("AS3").Value = 1250

# Row 4 updates
$ws.Range("AA4").Value = 1.91
$ws.Range("AB4").Value = 1.8
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 13
$ws.Range("AK4").Value = 15
$ws.Range("AO4").Value = 10
